$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "25.904.92"
Set-TextValue "E2" "  +0.27%  "
Set-TextValue "D3" "1.646.14"
Set-TextValue "E3" "  +0.72%  "
Set-TextValue "D4" "1.008"
Set-TextValue "E4" "  +0.54%  "
Set-TextValue "D5" "215.28"
Set-TextValue "E5" "  +0.01%  "
Set-TextValue "D6" "0.5081"
Set-TextValue "E6" "  +1.13%  "
Set-TextValue "D7" "1.006"
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "D9" "0.06412"
Set-TextValue "E9" "  +0.06%  "
Set-TextValue "D10" "19.74"
Set-TextValue "E10" "  +0.27%  "
Set-TextValue "D11" "0.07783"
Set-TextValue "E11" "  +1.34%  "
Set-TextValue "D12" "4.309"
Set-TextValue "E12" "  +1.53%  "
Set-TextValue "D13" "1.631.78"
Set-TextValue "E13" "  -0.49%  "
Set-TextValue "D14" "0.5465"
Set-TextValue "E14" "  +0.10%  "
Set-TextValue "D15" "0.0₅7894"
Set-TextValue "E15" "  -0.49%  "
Set-TextValue "D16" "65.19"
Set-TextValue "E16" "  +2.44%  "
Set-TextValue "D17" "25.983.32"
Set-TextValue "E18" "  +0.38%  "
Set-TextValue "D19" "197.19"
Set-TextValue "E19" "  -2.90%  "
Set-TextValue "D20" "4.414"
Set-TextValue "D21" "10.02"
Set-TextValue "E21" "  +0.81%  "
Set-TextValue "D22" "6.068"
Set-TextValue "E22" "  +1.53%  "
Set-TextValue "D23" "1.007"
Set-TextValue "E23" "  +0.37%  "
Set-TextValue "D24" "1.872"
Set-TextValue "E24" "  -3.05%  "
Set-TextValue "D25" "141.20"
Set-TextValue "E25" "  +0.02%  "
Set-TextValue "D26" "0.1145"
Set-TextValue "E26" "  +0.12%  "
Set-TextValue "D27" "6.901"
Set-TextValue "E27" "  +2.99%  "
Set-TextValue "D28" "15.75"
Set-TextValue "E28" "  +0.40%  "
Set-TextValue "D29" "0.05063"
Set-TextValue "E29" "  +1.41%  "
Set-TextValue "E30" "  +0.00%  "
Set-TextValue "D31" "3.272"
Set-TextValue "E31" "  +0.08%  "
Set-TextValue "D32" "3.209"
Set-TextValue "E32" "  +0.68%  "
Set-TextValue "E33" "  +0.41%  "
Set-TextValue "E34" "  +0.59%  "
Set-TextValue "D35" "0.8938"
Set-TextValue "D36" "2.602"
Set-TextValue "E36" "  -0.56%  "
Set-TextValue "B37" "Maker"
Set-TextValue "C37" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D37" "1.136.22"
Set-TextValue "E37" "  -3.05%  "
Set-TextValue "B38" "ImmutableX"
Set-TextValue "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.5542"
Set-TextValue "E38" "  -0.65%  "
Set-TextValue "D39" "0.01567"
Set-TextValue "E39" "  +0.45%  "
Set-TextValue "E40" "  +0.37%  "
Set-TextValue "D41" "5.672"
Set-TextValue "E41" "  +0.17%  "
Set-TextValue "D42" "0.8153"
Set-TextValue "E42" "  +1.07%  "
Set-TextValue "D43" "99.65"
Set-TextValue "E43" "  +0.26%  "
Set-TextValue "D44" "0.0₈123"
Set-TextValue "E44" "  +6.38%  "
Set-TextValue "D45" "1.782.66"
Set-TextValue "E45" "  +0.65%  "
Set-TextValue "D46" "0.4540"
Set-TextValue "E46" "  +0.60%  "
Set-TextValue "D47" "1.007"
Set-TextValue "E47" "  +0.02%  "
Set-TextValue "D48" "55.21"
Set-TextValue "E48" "  +0.56%  "
Set-TextValue "D49" "0.05092"
Set-TextValue "E49" "  +1.03%  "
Set-TextValue "D50" "1.008"
Set-TextValue "E50" "  +0.57%  "
Set-TextValue "E51" "  +3.25%  "
